$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This pharmacy report lists one product per row (columns B..N, merged),
# with column A holding a running sequence number and a totals/footer
# block right below the last data row. Two line items need to be removed
# from the table: "CETAL 250MG/5ML 60ML SUSP" and
# "OTRIVIN 0.05% PEDIATRIC NASAL DROPS 15 ML".

$namesToRemove = @(
    "CETAL 250MG/5ML 60ML SUSP",
    "OTRIVIN 0.05% PEDIATRIC NASAL DROPS 15 ML"
)

foreach ($name in $namesToRemove) {
    # Find the row whose product-name cell (column B) matches, then
    # delete the entire row - this shifts every following row (and the
    # totals/footer block) up by one.
    $found = $ws.Cells.Find($name)
    if ($found -ne $null) {
        $foundRow = $found.Row()
        $ws.Rows($foundRow).Delete()
    }
}

# Renumber the running sequence number in column A (1..8) for the rows
# that remain in the data block.
$dataRow = 4
$seq = 1
while ($true) {
    $nameVal = $ws.Cells.Item($dataRow, 2).Value()
    if ($nameVal -eq $null -or $nameVal -eq "") {
        break
    }
    $ws.Cells.Item($dataRow, 1).Value = $seq
    $seq = $seq + 1
    $dataRow = $dataRow + 1
}

# The totals row (column K, merged K:N) right below the data block must
# reflect the removal of the two rows' "سعر البيع" (price) values
# (31 + 24 = 55 less than before: 278.67 -> 223.67), and keeps its
# slightly taller row height.
$totalsRow = $dataRow
$ws.Cells.Item($totalsRow, 11).Value = 223.67
$ws.Rows($totalsRow).RowHeight = 26.25
